$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 93 (Indice 92): Mushuc Runa vs Tecnico U. ----
$ws.Cells.Item(93, 1).Value = 92
$ws.Cells.Item(93, 2).Value = "ecuador"
$ws.Cells.Item(93, 3).Value = "liga-pro"
$ws.Cells.Item(93, 4).Value = "'2023"
$ws.Cells.Item(93, 5).Value = 45235.77083333334
$ws.Cells.Item(93, 6).Value = "Mushuc Runa"
$ws.Cells.Item(93, 7).Value = 1
$ws.Cells.Item(93, 8).Value = "Tecnico U."
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 10).Value = 2.83
$ws.Cells.Item(93, 11).Value = "'29/10/2023 21:42"
$ws.Cells.Item(93, 12).Value = 3.2
$ws.Cells.Item(93, 13).Value = "'05/11/2023 18:21"
$ws.Cells.Item(93, 14).Value = 3.31
$ws.Cells.Item(93, 15).Value = "'29/10/2023 21:42"
$ws.Cells.Item(93, 16).Value = 3.33
$ws.Cells.Item(93, 17).Value = "'05/11/2023 18:21"
$ws.Cells.Item(93, 18).Value = 2.41
$ws.Cells.Item(93, 19).Value = "'29/10/2023 21:42"
$ws.Cells.Item(93, 20).Value = 2.31
$ws.Cells.Item(93, 21).Value = "'05/11/2023 18:24"
$ws.Cells.Item(93, 22).Value = "https://www.betexplorer.com/football/ecuador/liga-pro/mushuc-runa-tecnico-u/bwFtHQZT/"

# ---- Row 94 (Indice 93): Cumbaya vs Orense ----
$ws.Cells.Item(94, 1).Value = 93
$ws.Cells.Item(94, 2).Value = "ecuador"
$ws.Cells.Item(94, 3).Value = "liga-pro"
$ws.Cells.Item(94, 4).Value = "'2023"
$ws.Cells.Item(94, 5).Value = 45235.875
$ws.Cells.Item(94, 6).Value = "Cumbaya"
$ws.Cells.Item(94, 7).Value = 0
$ws.Cells.Item(94, 8).Value = "Orense"
$ws.Cells.Item(94, 9).Value = 1
$ws.Cells.Item(94, 10).Value = 3.34
$ws.Cells.Item(94, 11).Value = "'30/10/2023 00:12"
$ws.Cells.Item(94, 12).Value = 3.21
$ws.Cells.Item(94, 13).Value = "'05/11/2023 20:56"
$ws.Cells.Item(94, 14).Value = 2.99
$ws.Cells.Item(94, 15).Value = "'30/10/2023 00:12"
$ws.Cells.Item(94, 16).Value = 2.97
$ws.Cells.Item(94, 17).Value = "'05/11/2023 20:56"
$ws.Cells.Item(94, 18).Value = 2.28
$ws.Cells.Item(94, 19).Value = "'30/10/2023 00:12"
$ws.Cells.Item(94, 20).Value = 2.53
$ws.Cells.Item(94, 21).Value = "'05/11/2023 20:56"
$ws.Cells.Item(94, 22).Value = "https://www.betexplorer.com/football/ecuador/liga-pro/cumbaya-orense/EkkjyUKo/"

# Re-apply the formatting of the previous data row (92) onto the two new rows so they
# carry the same per-cell styles (bold/bordered/centered Indice in col A, date-time
# number format in col E, default style elsewhere) as the rest of the sheet. This also
# strips the transient "quote prefix" formatting that the leading apostrophes above
# applied to cells that needed to stay text (D, K, M, O, Q, S, U) instead of being
# auto-converted to numbers by the smart Value setter.
$ws.Range("A92:V92").Copy()
$ws.Range("A93:V94").PasteSpecial(-4122)  # xlPasteFormats
